# Restore functions of setting trip type, modify test_search_flight
#
# The input-data sheet previously had "extra" end_month/end_day values
# (columns I and J) hard-coded for rows 2-4, left over while the trip-type
# feature was disabled. Now that trip type handling is restored, those
# leftover one-way end-date values are no longer needed for those rows, so
# clear them back out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused end_month (I) / end_day (J) values on rows 2-4.
$ws.Range("I2:J2").ClearContents()
$ws.Range("I3:J3").ClearContents()
$ws.Range("I4:J4").ClearContents()

# Move the saved selection to J16 (matches the author's last selection).
[void]$ws.Range("J16").Select()
